$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$cellAddr,
        [string]$text
    )
    $range = $ws.Range($cellAddr)
    # Force a literal-text interpretation so ambiguous dd-mm-yyyy strings
    # (e.g. "01-08-2022") aren't auto-converted to date serials by Excel's
    # smart input parsing, then strip the temporary Text number format so
    # the cell's style/formatting is left exactly as it was before.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Update date strings from slash to dash format (rows 3..21, column A)
Set-TextValue "A3"  "28-07-2022"
Set-TextValue "A4"  "01-08-2022"
Set-TextValue "A5"  "04-08-2022"
Set-TextValue "A6"  "08-08-2022"
Set-TextValue "A7"  "11-08-2022"
Set-TextValue "A8"  "15-08-2022"
Set-TextValue "A9"  "18-08-2022"
Set-TextValue "A10" "22-08-2022"
Set-TextValue "A11" "25-08-2022"
Set-TextValue "A12" "29-08-2022"
Set-TextValue "A13" "01-09-2022"
Set-TextValue "A14" "05-09-2022"
Set-TextValue "A15" "08-09-2022"
Set-TextValue "A16" "12-09-2022"
Set-TextValue "A17" "15-09-2022"
Set-TextValue "A18" "19-09-2022"
Set-TextValue "A19" "22-09-2022"
Set-TextValue "A20" "26-09-2022"
Set-TextValue "A21" "29-09-2022"

# Update attendance counts for row 3 (Total Attendance Count and Invalid)
$ws.Range("D3").Value = 2
$ws.Range("G3").Value = 2
